# Atualização de bases das ligas, do dia: 24-02-2024 às 23:13
#
# The source data rows got re-sorted (e.g. by kickoff id) and several rows
# swapped places. Column A (row id), C (Div), D (Div Original Name) and
# E (Date) stay put; everything else (B id, F/G teams, H/I score, J result,
# K..AC odds) moves with the row it belongs to. This script reproduces that
# rotation by reading each row's current contents first (into an in-memory
# snapshot) and then writing the snapshot back out in the new row order, so
# within a rotation group no value is read after it has already been
# overwritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowSnapshot($sheet, [int]$row, [string[]]$cols) {
    $snap = @{}
    foreach ($c in $cols) {
        $snap[$c] = $sheet.Range("$c$row").Value()
    }
    return $snap
}

function Set-RowFromSnapshot($sheet, [int]$row, [string[]]$cols, $snap) {
    foreach ($c in $cols) {
        $sheet.Range("$c$row").Value = $snap[$c]
    }
}

# Columns present on completed fixtures (score + result recorded).
$fullCols = @('B','F','G','H','I','J','K','L','M','N','O','P','Q','R','S','T','U','V','W','X','Y','Z','AA','AB','AC')

# Columns present on not-yet-played fixtures (no FTHG/FTAG/FTR, no PL_AhOver/PL_AhUnder).
$partialCols = @('B','F','G','K','L','M','N','O','P','Q','R','S','T','U','V','W','X','Y','Z','AA')

# Each group is rotated: row -> gets the old contents of the row it maps to.
# Snapshot every row in the group before writing any of them back.

function Apply-Rotation($sheet, [int[]]$rows, [int[]]$sources, [string[]]$cols) {
    $snaps = @{}
    for ($i = 0; $i -lt $rows.Length; $i++) {
        $snaps[$rows[$i]] = Get-RowSnapshot $sheet $rows[$i] $cols
    }
    for ($i = 0; $i -lt $rows.Length; $i++) {
        Set-RowFromSnapshot $sheet $rows[$i] $cols $snaps[$sources[$i]]
    }
}

# Group 1: rows 586,587,588 -> 586 takes old 587, 587 takes old 588, 588 takes old 586
Apply-Rotation $ws @(586,587,588) @(587,588,586) $fullCols

# Group 2: rows 608,609 swap
Apply-Rotation $ws @(608,609) @(609,608) $fullCols

# Group 3: rows 613,614,617 -> 613 takes old 614, 614 takes old 617, 617 takes old 613
Apply-Rotation $ws @(613,614,617) @(614,617,613) $fullCols

# Group 4: rows 647,648 swap
Apply-Rotation $ws @(647,648) @(648,647) $fullCols

# Group 5: rows 649,650 swap
Apply-Rotation $ws @(649,650) @(650,649) $fullCols

# Group 6: rows 696,697 swap (future fixtures, no H/I/J/AB/AC columns)
Apply-Rotation $ws @(696,697) @(697,696) $partialCols

# Group 7: rows 700,701 swap (future fixtures, no H/I/J/AB/AC columns)
Apply-Rotation $ws @(700,701) @(701,700) $partialCols
